$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new data row above row 71 (shifts old rows 71-75 down
#    to 72-76) and clone the formatting/merges of the row that is
#    now below it (row 72, which used to be row 71) so the new row
#    looks exactly like the other product rows in the table.
# ------------------------------------------------------------------
$ws.Rows("71:71").Insert()

$ws.Range("A72:Q72").Copy()
$ws.Range("A71:Q71").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A71:B71").Merge()
$ws.Range("C71:G71").Merge()
$ws.Range("H71:K71").Merge()
$ws.Range("L71:M71").Merge()
$ws.Range("N71:O71").Merge()

# ------------------------------------------------------------------
# 2. Fill in the new product row (صوفى طويل جدا جدا).
# ------------------------------------------------------------------
$ws.Range("A71").Value = 65
$ws.Range("C71").Value = "صوفى طويل جدا جدا"
$ws.Range("H71").Value = "8:0"
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = "55.00"
$ws.Range("P71").Value = "55.0000"
$ws.Range("Q71").Value = 12

# ------------------------------------------------------------------
# 3. Update the grand-total cell (previously P74, now shifted to
#    P75) to reflect the newly added row.
# ------------------------------------------------------------------
$ws.Range("P75").Value = 3785.625

# ------------------------------------------------------------------
# 4. Refresh the generated-on timestamp in the footer (previously
#    row 75, now shifted to row 76).
# ------------------------------------------------------------------
$ws.Range("A76").Value = "Tuesday, 19 August, 2025 6:30 PM"
